# Vælg tidspunkt fra listbox.
# Adds a new data row (row 15) to the sensor data sheet, mirroring a row
# picked/entered via the workbook's time-selection listbox.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 15

$ws.Cells.Item($row, 1).Value = 5
$ws.Cells.Item($row, 2).Value = "Luna"

# Tidspunkt - date/time value, reuse the same style as the cells above (m/d/yyyy h:mm)
$ws.Range("C14").Copy()
$ws.Cells.Item($row, 3).PasteSpecial(-4122)
$ws.Cells.Item($row, 3).Value = 43382.772916666669

$ws.Cells.Item($row, 4).Value = 0
$ws.Cells.Item($row, 5).Value = 0
$ws.Cells.Item($row, 6).Value = 1
$ws.Cells.Item($row, 7).Value = 0

# Varighed for arbejdsgang - stored as literal text "10.23" (not a numeric time)
$ws.Range("H14").Copy()
$ws.Cells.Item($row, 8).PasteSpecial(-4122)
$ws.Cells.Item($row, 8).Formula = "=""10.23"""
$ws.Cells.Item($row, 8).Copy()
$ws.Cells.Item($row, 8).PasteSpecial(-4163)

$ws.Cells.Item($row, 9).Value = "Toiletbesøg"
$ws.Cells.Item($row, 10).Value = 1

# Tid med borger - stored as literal text "00.40"
$ws.Cells.Item($row, 11).Formula = "=""00.40"""
$ws.Cells.Item($row, 11).Copy()
$ws.Cells.Item($row, 11).PasteSpecial(-4163)

$ws.Cells.Item($row, 12).Value = "Alm"

$excel.CutCopyMode = 0
[void]$ws.Range("L16").Select()
